# Add three more failed-login attempts ("Balu") to the invalidLoginTest
# sheet (rows 4-6) and move the active-tab/selection from validLoginTest
# (sheet2) over to invalidLoginTest (sheet1).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New data rows: Username | Password | expected result
$ws1.Range("A4").Value = "Balu"
$ws1.Range("B4").Value = "balu123"
$ws1.Range("C4").Value = "Invalid credentials"

$ws1.Range("A5").Value = "Balu"
$ws1.Range("B5").Value = "balu124"
$ws1.Range("C5").Value = "Invalid credentials"

$ws1.Range("A6").Value = "Balu"
$ws1.Range("B6").Value = "balu125"
$ws1.Range("C6").Value = "Invalid credentials"

# Make invalidLoginTest the active sheet/tab and select C3:C6 there
# (this also clears tabSelected on the previously-active sheet2).
$ws1.Activate()
$ws1.Range("C3:C6").Select()
